# Daily "cryptos" price/volume refresh (GitHub Actions scraper update).
# Updates the Price (D) and Volume(1h) (E) columns for each coin row, and
# swaps the Maker / RenderToken rows (42/43) which changed rank order.
#
# Note: several Price values look like plain numbers (e.g. "290.78",
# "0.1000") but must stay literal text (matching the source inlineStr
# cells, including trailing zeros). Excel's Range.Value setter normally
# "smart types" such strings into real numbers, so for those cells we
# force the Text number format first, exactly as you would in the UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.116.01'
$ws.Range('E2').Value = '  +1.24%  '

$ws.Range('D3').Value = '2.219.37'
$ws.Range('E3').Value = '  +0.17%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '290.78'
$ws.Range('E5').Value = '  -2.53%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '88.29'
$ws.Range('E6').Value = '  +5.54%  '

$ws.Range('E7').Value = '  +0.25%  '

$ws.Range('E8').Value = '  -0.09%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.472'
$ws.Range('E9').Value = '  +1.65%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '30.82'
$ws.Range('E10').Value = '  +4.09%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0784'
$ws.Range('E11').Value = '  +0.41%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.92'
$ws.Range('E12').Value = '  +4.16%  '

$ws.Range('E13').Value = '  +2.69%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.48'
$ws.Range('E14').Value = '  +3.51%  '

$ws.Range('D15').Value = '2.556.63'
$ws.Range('E15').Value = '  -0.24%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.05'
$ws.Range('E16').Value = '  -0.37%  '

$ws.Range('D17').Value = '2.212.75'
$ws.Range('E17').Value = '  +0.20%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.730'
$ws.Range('E18').Value = '  +1.86%  '

$ws.Range('D19').Value = '40.024.71'
$ws.Range('E19').Value = '  +1.16%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.95'
$ws.Range('E20').Value = '  +15.19%  '

$ws.Range('D21').Value = '0.0₃0889'
$ws.Range('E21').Value = '  +1.27%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.83'
$ws.Range('E22').Value = '  +1.60%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.73'
$ws.Range('E23').Value = '  +1.14%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '235.87'
$ws.Range('E24').Value = '  +1.41%  '

$ws.Range('E25').Value = '  +0.09%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.47'
$ws.Range('E26').Value = '  +1.95%  '

$ws.Range('E27').Value = '  +1.08%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '22.69'
$ws.Range('E28').Value = '  -0.10%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.20'
$ws.Range('E29').Value = '  +1.28%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.25'
$ws.Range('E30').Value = '  +1.10%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '153.31'
$ws.Range('E31').Value = '  +2.64%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '32.31'
$ws.Range('E32').Value = '  +0.19%  '

$ws.Range('E33').Value = '  -0.09%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.97'
$ws.Range('E34').Value = '  +2.94%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0721'
$ws.Range('E35').Value = '  +2.89%  '

$ws.Range('E36').Value = '  -0.09%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.85'
$ws.Range('E37').Value = '  +7.82%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '16.07'
$ws.Range('E38').Value = '  -0.60%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.112'
$ws.Range('E39').Value = '  +0.81%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.1000'
$ws.Range('E40').Value = '  +3.02%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.71'
$ws.Range('E41').Value = '  +3.62%  '

# Row 42/43 swap: RenderToken overtook Maker in rank.
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.84'
$ws.Range('E42').Value = '  +5.19%  '

$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.095.97'
$ws.Range('E43').Value = '  +8.70%  '

$ws.Range('E44').Value = '  +2.58%  '

$ws.Range('E45').Value = '  +1.61%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.89'
$ws.Range('E46').Value = '  +7.35%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '17.64'
$ws.Range('E47').Value = '  +8.16%  '

$ws.Range('E48').Value = '  +3.00%  '

$ws.Range('D49').Value = '2.428.37'
$ws.Range('E49').Value = '  -0.34%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '69.75'
$ws.Range('E50').Value = '  -1.46%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '88.96'
$ws.Range('E51').Value = '  +0.48%  '
